$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.157994747161865
$ws.Range("B1").Value = 2.417201042175293
$ws.Range("C1").Value = 2.527528285980225
$ws.Range("D1").Value = 3.228439807891846
$ws.Range("E1").Value = 2.321007251739502
